# Apply the "symbol list" refresh for Mon Dec 19 22:07:11 UTC 2022 (GitHub Actions scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and G (Hora) hold digit-only text in the source data; force the
# Text number format on the cells being touched so Excel keeps the new values stored as
# text instead of auto-coercing them to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D4:D6").NumberFormat = "@"
$ws.Range("D8:D27").NumberFormat = "@"
$ws.Range("D40:D45").NumberFormat = "@"
$ws.Range("D47:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# "Hora" (hour) column bumps from 21 to 22 for every data row.
$ws.Range("G2:G51").Value = "22"

# Row 2
$ws.Range("D2").Value = "242.91"

# Row 4
$ws.Range("D4").Value = "5.201"

# Row 5
$ws.Range("D5").Value = "0.05600"

# Row 6
$ws.Range("D6").Value = "3.371"

# Row 8
$ws.Range("D8").Value = "0.8056"

# Row 9
$ws.Range("D9").Value = "0.9704"

# Row 10
$ws.Range("D10").Value = "0.1405"

# Row 11
$ws.Range("D11").Value = "0.07329"

# Row 12
$ws.Range("D12").Value = "0.03121"

# Row 13
$ws.Range("B13").Value = "ProBitToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D13").Value = "0.1287"
$ws.Range("E13").Value = "12ProBitTokenPROB"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03048"
$ws.Range("E14").Value = "13BitrueCoinBTR"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09288"
$ws.Range("E15").Value = "14BitMartTokenBMX"

# Row 16
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.610"
$ws.Range("E16").Value = "15MCDexMCB"

# Row 17
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001642"
$ws.Range("E17").Value = "16BitForexTokenBF"

# Row 18
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04710"
$ws.Range("E18").Value = "17CoinExTokenCET"

# Row 19
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value = "0.0005753"
$ws.Range("E19").Value = "18OneONE"

# Row 20
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value = "0.006401"
$ws.Range("E20").Value = "19TigerCashTCH"

# Row 21
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "0.004982"
$ws.Range("E21").Value = "20HotbitTokenHTB"

# Row 22
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "0.001042"
$ws.Range("E22").Value = "21BitKanKAN"

# Row 23
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "0.0001500"
$ws.Range("E23").Value = "22NitroExNTX"

# Row 24
$ws.Range("B24").Value = "UpBots"
$ws.Range("C24").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D24").Value = "0.0003102"
$ws.Range("E24").Value = "23UpBotsUBXT"

# Row 25
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "3.754"
$ws.Range("E25").Value = "24LEOLEO"

# Row 26
$ws.Range("B26").Value = "BTSEToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D26").Value = "2.097"
$ws.Range("E26").Value = "25BTSETokenBTSE"

# Row 27
$ws.Range("B27").Value = "BitpandaEcosystemToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D27").Value = "0.3261"
$ws.Range("E27").Value = "26BitpandaEcosystemTokenBEST"

# Row 40
$ws.Range("D40").Value = "0.03888"

# Row 41
$ws.Range("D41").Value = "0.006892"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1037"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002900"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44
$ws.Range("D44").Value = "0.008410"

# Row 45
$ws.Range("D45").Value = "0.00005934"

# Row 47
$ws.Range("D47").Value = "0.0005503"

# Row 48
$ws.Range("D48").Value = "0.6828"

# Row 49
$ws.Range("D49").Value = "0.09376"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"

# Row 50
$ws.Range("D50").Value = "0.00002101"

# Row 51
$ws.Range("D51").Value = "0.01011"
